$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138 (shifts rows 138:239 down to 139:240,
# extends the used range / dimension to A1:R240).
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new record's data.
$ws.Cells.Item(138, 1).Value = 8
$ws.Cells.Item(138, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 45072
$ws.Cells.Item(138, 5).Value = 4
$ws.Cells.Item(138, 6).Value = 100112001
$ws.Cells.Item(138, 7).Value = "Berenjena"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 480
$ws.Cells.Item(138, 11).Value = 9500
$ws.Cells.Item(138, 12).Value = 10000
$ws.Cells.Item(138, 13).Value = 9750
$ws.Cells.Item(138, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(138, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(138, 16).Value = 195
$ws.Cells.Item(138, 17).Value = 50
$ws.Cells.Item(138, 18).Value = "Hortaliza"
